# Daily attendance processing - 2026-01-06 21:59:01
# Rotate the "Recorded By" (column G) comma-separated list so that the
# last entry moves to the front, for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
            $cell.Value2 = [string]::Join(", ", $rotated)
        }
    }
}
